$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: new entry for 14.01.2020
$ws.Range("A15").Value = "Di"
$ws.Range("B15").Value = "14.01.2020"
$ws.Range("C15").Value = 0.32291666666666669
$ws.Range("D15").Value = 0.54861111111111105
$ws.Range("E15").Value = "Sprintplanung (Akzeptanzkriterien, Flip Chart, Story Points,…), Start des 3. Sprints, Programmierarbeiten an Story"

$ws.Rows.Item(15).RowHeight = 27

# Row 16: new entry for 15.01.2020
$ws.Range("A16").Value = 43845
$ws.Range("B16").Value = "15.01.2020"
$ws.Range("C16").Value = 0.40277777777777773
$ws.Range("D16").Value = 0.47222222222222227
$ws.Range("E16").Value = "Programmieren an Story - Admin User anzeigen, bearbeiten"

# Update current selection to match the saved view state
$ws.Range("E17").Select()
